$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'275.32"
$ws.Range("E2").Value = "'-0.76%"
$ws.Range("G2").Value = "'8"
$ws.Range("D3").Value = "'26.67"
$ws.Range("E3").Value = "'-2.34%"
$ws.Range("G3").Value = "'8"
$ws.Range("D4").Value = "'4.879"
$ws.Range("E4").Value = "'1.77%"
$ws.Range("G4").Value = "'8"
$ws.Range("D5").Value = "'0.06331"
$ws.Range("E5").Value = "'1.23%"
$ws.Range("G5").Value = "'8"
$ws.Range("D6").Value = "'6.881"
$ws.Range("E6").Value = "'-0.61%"
$ws.Range("G6").Value = "'8"
$ws.Range("D7").Value = "'3.318"
$ws.Range("E7").Value = "'1.60%"
$ws.Range("G7").Value = "'8"
$ws.Range("D8").Value = "'1.253"
$ws.Range("E8").Value = "'32.87%"
$ws.Range("G8").Value = "'8"
$ws.Range("D9").Value = "'0.8698"
$ws.Range("E9").Value = "'-1.24%"
$ws.Range("G9").Value = "'8"
$ws.Range("D10").Value = "'0.1590"
$ws.Range("E10").Value = "'8.86%"
$ws.Range("G10").Value = "'8"
$ws.Range("D11").Value = "'0.05199"
$ws.Range("E11").Value = "'-2.15%"
$ws.Range("G11").Value = "'8"
$ws.Range("D12").Value = "'0.07439"
$ws.Range("E12").Value = "'1.12%"
$ws.Range("G12").Value = "'8"
$ws.Range("D13").Value = "'0.02927"
$ws.Range("E13").Value = "'-5.92%"
$ws.Range("G13").Value = "'8"
$ws.Range("D14").Value = "'0.09039"
$ws.Range("E14").Value = "'-0.22%"
$ws.Range("G14").Value = "'8"
$ws.Range("D15").Value = "'0.001589"
$ws.Range("E15").Value = "'2.32%"
$ws.Range("G15").Value = "'8"
$ws.Range("D16").Value = "'0.0006338"
$ws.Range("E16").Value = "'0.95%"
$ws.Range("G16").Value = "'8"
$ws.Range("D17").Value = "'0.006015"
$ws.Range("E17").Value = "'5.16%"
$ws.Range("G17").Value = "'8"
$ws.Range("D18").Value = "'3.451"
$ws.Range("E18").Value = "'0.32%"
$ws.Range("G18").Value = "'8"
$ws.Range("D19").Value = "'2.272"
$ws.Range("E19").Value = "'0.35%"
$ws.Range("G19").Value = "'8"
$ws.Range("D20").Value = "'0.3113"
$ws.Range("E20").Value = "'-1.05%"
$ws.Range("G20").Value = "'8"
$ws.Range("D21").Value = "'0.1336"
$ws.Range("E21").Value = "'1.85%"
$ws.Range("G21").Value = "'8"
$ws.Range("D22").Value = "'3.906"
$ws.Range("E22").Value = "'1.49%"
$ws.Range("G22").Value = "'8"
$ws.Range("D23").Value = "'0.04359"
$ws.Range("E23").Value = "'0.57%"
$ws.Range("G23").Value = "'8"
$ws.Range("D24").Value = "'0.001180"
$ws.Range("E24").Value = "'0.17%"
$ws.Range("G24").Value = "'8"
$ws.Range("D25").Value = "'0.004248"
$ws.Range("E25").Value = "'-0.68%"
$ws.Range("G25").Value = "'8"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("G26").Value = "'8"
$ws.Range("E27").Value = "'-0.25%"
$ws.Range("G27").Value = "'8"
$ws.Range("G28").Value = "'8"
$ws.Range("G29").Value = "'8"
$ws.Range("G30").Value = "'8"
$ws.Range("G31").Value = "'8"
$ws.Range("G32").Value = "'8"
$ws.Range("G33").Value = "'8"
$ws.Range("G34").Value = "'8"
$ws.Range("G35").Value = "'8"
$ws.Range("G36").Value = "'8"
$ws.Range("G37").Value = "'8"
$ws.Range("G38").Value = "'8"
$ws.Range("G39").Value = "'8"
$ws.Range("D40").Value = "'0.04108"
$ws.Range("E40").Value = "'1.47%"
$ws.Range("G40").Value = "'8"
$ws.Range("D41").Value = "'0.006739"
$ws.Range("E41").Value = "'0.32%"
$ws.Range("G41").Value = "'8"
$ws.Range("D42").Value = "'0.1164"
$ws.Range("E42").Value = "'0.85%"
$ws.Range("G42").Value = "'8"
$ws.Range("E43").Value = "'0.76%"
$ws.Range("G43").Value = "'8"
$ws.Range("D44").Value = "'0.01076"
$ws.Range("E44").Value = "'-17.25%"
$ws.Range("G44").Value = "'8"
$ws.Range("D45").Value = "'0.00005288"
$ws.Range("E45").Value = "'3.53%"
$ws.Range("G45").Value = "'8"
$ws.Range("D46").Value = "'0.02103"
$ws.Range("E46").Value = "'-29.60%"
$ws.Range("G46").Value = "'8"
$ws.Range("E47").Value = "'-37.33%"
$ws.Range("G47").Value = "'8"
$ws.Range("G48").Value = "'8"
$ws.Range("G49").Value = "'8"
$ws.Range("G50").Value = "'8"
$ws.Range("G51").Value = "'8"
